$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns AD/AE/AF -> Wins/Losses/Ties
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the existing bold/bordered/centered header style used by A1:AC1
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows 2-44: team record values (Wins=87, Losses=76, Ties=0)
$ws.Range("AD2:AD44").Value = 87
$ws.Range("AE2:AE44").Value = 76
$ws.Range("AF2:AF44").Value = 0
